# Rename the worksheet from "Sheet1" to "Analysis" to better reflect its
# content (non-conformities / threat-assessment analysis).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Analysis"
